$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for "Feria Lagunitas de Puerto
# Montt - Cilantro". It belongs at the top of the data block (row 81, right
# after the header + region rows above it), so insert a fresh row there and
# push the existing data (old rows 81-175) down by one — the row that used
# to be last (old row 175) ends up as the new row 176.
$ws.Rows("81").Insert()

# Populate the newly inserted row 81 with the new observation.
$ws.Range("A81").Value = 4
$ws.Range("B81").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C81").Value = "Los Lagos"
$ws.Range("D81").Value = 44483
$ws.Range("E81").Value = 10
$ws.Range("F81").Value = 100112040
$ws.Range("G81").Value = "Cilantro"
$ws.Range("H81").Value = "Sin especificar"
$ws.Range("I81").Value = "Primera"
$ws.Range("J81").Value = 150
$ws.Range("K81").Value = 10000
$ws.Range("L81").Value = 10000
$ws.Range("M81").Value = 10000
$ws.Range("N81").Value = "`$/caja 36 atados"
$ws.Range("O81").Value = "Región Metropolitana"
$ws.Range("P81").Value = 278
$ws.Range("Q81").Value = 36
$ws.Range("R81").Value = "Hortaliza"
